$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2832150.2
$ws.Range("I74").Value = 3185694
$ws.Range("J74").Value = 3800
$ws.Range("K74").Value = 3185694
$ws.Range("L74").Value = 3800
$ws.Range("M74").Value = -3184758
$ws.Range("N74").Value = -5672
$ws.Range("H77").Value = 2832150.2
$ws.Range("I77").Value = 3185694
$ws.Range("J77").Value = 3800
$ws.Range("K77").Value = 15928470
$ws.Range("L77").Value = 19000
$ws.Range("M77").Value = -15923790
$ws.Range("N77").Value = -28360
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
$ws.Range("H132").Value = 1899.8108
$ws.Range("I132").Value = 1957.5938
$ws.Range("K132").Value = 5872.7814
$ws.Range("M132").Value = -3342.7814
$ws.Range("H137").Value = 1072.738
$ws.Range("I137").Value = 892.25
$ws.Range("J137").Value = 1650.3
$ws.Range("K137").Value = 2676.75
$ws.Range("L137").Value = 4950.9
$ws.Range("M137").Value = -126.75
$ws.Range("N137").Value = -10050.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5485640
$ws.Range("I32").Value = 6083313
$ws.Range("J32").Value = 6966.6665
$ws.Range("K32").Value = 6083313
$ws.Range("L32").Value = 6966.6665
$ws.Range("M32").Value = -6083026
$ws.Range("N32").Value = -7540.6665
$ws.Range("H74").Value = 769.875
$ws.Range("I74").Value = 769.10254
$ws.Range("K74").Value = 769.10254
$ws.Range("M74").Value = 104.89746
$ws.Range("H77").Value = 769.875
$ws.Range("I77").Value = 769.10254
$ws.Range("K77").Value = 3845.5127
$ws.Range("M77").Value = 522.4873000000002
$ws.Range("H102").Value = 2458.158
$ws.Range("I102").Value = 2423.2354
$ws.Range("J102").Value = 2755
$ws.Range("K102").Value = 2423.2354
$ws.Range("L102").Value = 2755
$ws.Range("M102").Value = -801.2354
$ws.Range("N102").Value = -5999
$ws.Range("H132").Value = 1330.58
$ws.Range("I132").Value = 1047.075
$ws.Range("J132").Value = 2464.6
$ws.Range("K132").Value = 3141.225
$ws.Range("L132").Value = 7393.799999999999
$ws.Range("M132").Value = -611.2250000000004
$ws.Range("N132").Value = -12453.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10454.625
$ws.Range("I105").Value = 10000
$ws.Range("J105").Value = 11363.875
$ws.Range("K105").Value = 10000
$ws.Range("L105").Value = 11363.875
$ws.Range("M105").Value = -8253
$ws.Range("N105").Value = -14857.875
$ws.Range("H132").Value = 151000
$ws.Range("J132").Value = 151000
$ws.Range("L132").Value = 151000
$ws.Range("N132").Value = -161120
$ws.Range("H134").Value = 22051.633
$ws.Range("I134").Value = 1555.5853
$ws.Range("J134").Value = 127093.875
$ws.Range("K134").Value = 4666.7559
$ws.Range("L134").Value = 381281.625
$ws.Range("M134").Value = -2131.7559
$ws.Range("N134").Value = -386351.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 373.8
$ws.Range("I22").Value = 354.75
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 354.75
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -4.75
$ws.Range("N22").Value = -1150
$ws.Range("H31").Value = 2796.318
$ws.Range("I31").Value = 2810.923
$ws.Range("J31").Value = 2682.4
$ws.Range("K31").Value = 2810.923
$ws.Range("L31").Value = 2682.4
$ws.Range("M31").Value = -2515.923
$ws.Range("N31").Value = -3272.4
$ws.Range("H34").Value = 2796.318
$ws.Range("I34").Value = 2810.923
$ws.Range("J34").Value = 2682.4
$ws.Range("K34").Value = 2810.923
$ws.Range("L34").Value = 2682.4
$ws.Range("M34").Value = -2608.923
$ws.Range("N34").Value = -3086.4
$ws.Range("H58").Value = 3152.3901
$ws.Range("I58").Value = 623.42426
$ws.Range("J58").Value = 13584.375
$ws.Range("K58").Value = 623.42426
$ws.Range("L58").Value = 13584.375
$ws.Range("M58").Value = -420.42426
$ws.Range("N58").Value = -13990.375
$ws.Range("H132").Value = 3652.2
$ws.Range("I132").Value = 2006
$ws.Range("J132").Value = 4749.6665
$ws.Range("K132").Value = 6018
$ws.Range("L132").Value = 14248.9995
$ws.Range("M132").Value = -3488
$ws.Range("N132").Value = -19308.9995
$ws.Range("H134").Value = 6222.3335
$ws.Range("I134").Value = 3667
$ws.Range("J134").Value = 11333
$ws.Range("K134").Value = 11001
$ws.Range("L134").Value = 33999
$ws.Range("M134").Value = -8466
$ws.Range("N134").Value = -39069
$ws.Range("H136").Value = 3152.3901
$ws.Range("I136").Value = 623.42426
$ws.Range("J136").Value = 13584.375
$ws.Range("K136").Value = 1870.27278
$ws.Range("L136").Value = 40753.125
$ws.Range("M136").Value = 679.72722
$ws.Range("N136").Value = -45853.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 73408.28999999999
$ws.Range("I109").Value = 111635.11
$ws.Range("J109").Value = 4600
$ws.Range("K109").Value = 334905.33
$ws.Range("L109").Value = 13800
$ws.Range("M109").Value = -333865.33
$ws.Range("N109").Value = -15880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4584
$ws.Range("I70").Value = 4490
$ws.Range("J70").Value = 4678
$ws.Range("K70").Value = 4490
$ws.Range("L70").Value = 4678
$ws.Range("M70").Value = -4220
$ws.Range("N70").Value = -5218
$ws.Range("H73").Value = 4584
$ws.Range("I73").Value = 4490
$ws.Range("J73").Value = 4678
$ws.Range("K73").Value = 4490
$ws.Range("L73").Value = 4678
$ws.Range("M73").Value = -3554
$ws.Range("N73").Value = -6550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1652.6842
$ws.Range("I16").Value = 1085.7858
$ws.Range("K16").Value = 1085.7858
$ws.Range("M16").Value = -915.7858000000001
$ws.Range("H93").Value = 970.4583
$ws.Range("I93").Value = 940.3182
$ws.Range("J93").Value = 1302
$ws.Range("K93").Value = 940.3182
$ws.Range("L93").Value = 1302
$ws.Range("M93").Value = 307.6818
$ws.Range("N93").Value = -3798
$ws.Range("H136").Value = 1915.6666
$ws.Range("I136").Value = 1355.3658
$ws.Range("J136").Value = 3351.4375
$ws.Range("K136").Value = 4066.0974
$ws.Range("L136").Value = 10054.3125
$ws.Range("M136").Value = -1516.0974
$ws.Range("N136").Value = -15154.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 825
$ws.Range("I126").Value = 766.6667
$ws.Range("J126").Value = 860
$ws.Range("K126").Value = 2300.0001
$ws.Range("L126").Value = 2580
$ws.Range("M126").Value = 169.9998999999998
$ws.Range("N126").Value = -7520
$ws.Range("H132").Value = 1106.3877
$ws.Range("I132").Value = 758.3721
$ws.Range("J132").Value = 3600.5
$ws.Range("K132").Value = 2275.1163
$ws.Range("L132").Value = 10801.5
$ws.Range("M132").Value = 254.8836999999999
$ws.Range("N132").Value = -15861.5
$ws.Range("H136").Value = 1168.8363
$ws.Range("I136").Value = 1037.159
$ws.Range("J136").Value = 1695.5454
$ws.Range("K136").Value = 3111.477
$ws.Range("L136").Value = 5086.6362
$ws.Range("M136").Value = -561.4770000000003
$ws.Range("N136").Value = -10186.6362
$ws.Range("H140").Value = 44034.168
$ws.Range("J140").Value = 44034.168
$ws.Range("L140").Value = 44034.168
$ws.Range("N140").Value = -54394.168
